{"js": "// Replace each old equation text with its corresponding new equation text.\n// Each \"NNN\u00d7N=\" string is unique within the document, so a simple\n// search-and-replace (matching the whole substring, case-sensitive) is safe.\nconst replacements = [\n  [\"531\u00d74=\", \"698\u00d75=\"],\n  [\"586\u00d74=\", \"521\u00d72=\"],\n  [\"787\u00d73=\", \"337\u00d78=\"],\n  [\"178\u00d72=\", \"779\u00d77=\"],\n  [\"947\u00d78=\", \"319\u00d76=\"],\n  [\"361\u00d79=\", \"855\u00d79=\"],\n  [\"250\u00d76=\", \"506\u00d76=\"],\n  [\"833\u00d79=\", \"134\u00d77=\"],\n  [\"802\u00d76=\", \"120\u00d75=\"],\n  [\"619\u00d76=\", \"700\u00d76=\"],\n  [\"484\u00d78=\", \"743\u00d76=\"],\n  [\"125\u00d79=\", \"734\u00d78=\"],\n  [\"321\u00d77=\", \"344\u00d73=\"],\n  [\"691\u00d72=\", \"634\u00d77=\"],\n  [\"543\u00d77=\", \"663\u00d74=\"],\n  [\"781\u00d77=\", \"148\u00d79=\"],\n  [\"509\u00d72=\", \"560\u00d78=\"],\n  [\"568\u00d73=\", \"134\u00d78=\"],\n  [\"667\u00d73=\", \"327\u00d74=\"],\n  [\"417\u00d78=\", \"542\u00d75=\"],\n  [\"766\u00d72=\", \"480\u00d74=\"],\n  [\"572\u00d79=\", \"292\u00d79=\"],\n  [\"402\u00d77=\", \"295\u00d72=\"],\n  [\"157\u00d77=\", \"807\u00d79=\"],\n  [\"632\u00d75=\", \"876\u00d74=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const found = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each old equation text with its corresponding new equation text.\n# Each \"NNN\u00d7N=\" string is unique within the document, so a simple\n# Find/Replace (whole match, case-sensitive) is safe for every pair.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"531\u00d74=\", \"698\u00d75=\"),\n    @(\"586\u00d74=\", \"521\u00d72=\"),\n    @(\"787\u00d73=\", \"337\u00d78=\"),\n    @(\"178\u00d72=\", \"779\u00d77=\"),\n    @(\"947\u00d78=\", \"319\u00d76=\"),\n    @(\"361\u00d79=\", \"855\u00d79=\"),\n    @(\"250\u00d76=\", \"506\u00d76=\"),\n    @(\"833\u00d79=\", \"134\u00d77=\"),\n    @(\"802\u00d76=\", \"120\u00d75=\"),\n    @(\"619\u00d76=\", \"700\u00d76=\"),\n    @(\"484\u00d78=\", \"743\u00d76=\"),\n    @(\"125\u00d79=\", \"734\u00d78=\"),\n    @(\"321\u00d77=\", \"344\u00d73=\"),\n    @(\"691\u00d72=\", \"634\u00d77=\"),\n    @(\"543\u00d77=\", \"663\u00d74=\"),\n    @(\"781\u00d77=\", \"148\u00d79=\"),\n    @(\"509\u00d72=\", \"560\u00d78=\"),\n    @(\"568\u00d73=\", \"134\u00d78=\"),\n    @(\"667\u00d73=\", \"327\u00d74=\"),\n    @(\"417\u00d78=\", \"542\u00d75=\"),\n    @(\"766\u00d72=\", \"480\u00d74=\"),\n    @(\"572\u00d79=\", \"292\u00d79=\"),\n    @(\"402\u00d77=\", \"295\u00d72=\"),\n    @(\"157\u00d77=\", \"807\u00d79=\"),\n    @(\"632\u00d75=\", \"876\u00d74=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
